$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '26.816.58'
Set-TextValue 'E2' '  -3.05%  '
Set-TextValue 'D3' '1.857.24'
Set-TextValue 'E3' '  -2.04%  '
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  +0.13%  '
Set-TextValue 'D5' '305.06'
Set-TextValue 'E5' '  -1.71%  '
Set-TextValue 'E6' '  +0.13%  '
Set-TextValue 'D7' '0.5079'
Set-TextValue 'E7' '  -3.46%  '
Set-TextValue 'D8' '0.3653'
Set-TextValue 'E8' '  -3.71%  '
Set-TextValue 'E9' '  -1.44%  '
Set-TextValue 'B10' 'Solana'
Set-TextValue 'C10' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D10' '20.74'
Set-TextValue 'E10' '  -1.52%  '
Set-TextValue 'B11' 'Polygon'
Set-TextValue 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D11' '0.8874'
Set-TextValue 'E11' '  -1.46%  '
Set-TextValue 'D12' '1.867.85'
Set-TextValue 'E12' '  -1.55%  '
Set-TextValue 'D13' '0.07497'
Set-TextValue 'E13' '  -1.68%  '
Set-TextValue 'D14' '5.238'
Set-TextValue 'D15' '90.49'
Set-TextValue 'E15' '  -1.27%  '
Set-TextValue 'E16' '  +0.17%  '
Set-TextValue 'D17' '0.000008528'
Set-TextValue 'E17' '  -1.53%  '
Set-TextValue 'E18' '  -1.88%  '
Set-TextValue 'D19' '1.001'
Set-TextValue 'E19' '  +0.15%  '
Set-TextValue 'D20' '26.859.05'
Set-TextValue 'E20' '  -3.03%  '
Set-TextValue 'D21' '4.995'
Set-TextValue 'E21' '  -2.84%  '
Set-TextValue 'D22' '2.104.37'
Set-TextValue 'E22' '  -1.07%  '
Set-TextValue 'D23' '10.28'
Set-TextValue 'E23' '  -4.85%  '
Set-TextValue 'D24' '6.441'
Set-TextValue 'E24' '  -2.37%  '
Set-TextValue 'D25' '1.817'
Set-TextValue 'E25' '  -2.20%  '
Set-TextValue 'D26' '145.79'
Set-TextValue 'E26' '  -4.75%  '
Set-TextValue 'D27' '17.82'
Set-TextValue 'E27' '  -2.30%  '
Set-TextValue 'D28' '2.047'
Set-TextValue 'E28' '  -5.61%  '
Set-TextValue 'D29' '112.66'
Set-TextValue 'E29' '  -1.18%  '
Set-TextValue 'D30' '4.619'
Set-TextValue 'E30' '  -4.03%  '
Set-TextValue 'D31' '4.663'
Set-TextValue 'E31' '  -3.17%  '
Set-TextValue 'D32' '0.09210'
Set-TextValue 'E32' '  +0.66%  '
Set-TextValue 'D33' '0.05106'
Set-TextValue 'E33' '  -3.03%  '
Set-TextValue 'D34' '3.064'
Set-TextValue 'E34' '  -2.90%  '
Set-TextValue 'D35' '1.148'
Set-TextValue 'E35' '  -5.93%  '
Set-TextValue 'D36' '0.7314'
Set-TextValue 'E36' '  -5.15%  '
Set-TextValue 'E37' '  +3.72%  '
Set-TextValue 'D38' '0.02014'
Set-TextValue 'E38' '  -3.17%  '
Set-TextValue 'E39' '  -3.77%  '
Set-TextValue 'D40' '1.073'
Set-TextValue 'E40' '  -1.44%  '
Set-TextValue 'D41' '0.5309'
Set-TextValue 'E41' '  -4.47%  '
Set-TextValue 'D42' '118.01'
Set-TextValue 'E42' '  +0.61%  '
Set-TextValue 'D43' '6.461'
Set-TextValue 'E43' '  -3.32%  '
Set-TextValue 'D44' '8.389'
Set-TextValue 'E44' '  -3.62%  '
Set-TextValue 'D45' '0.1469'
Set-TextValue 'E45' '  -2.58%  '
Set-TextValue 'E46' '  +0.24%  '
Set-TextValue 'D47' '0.4626'
Set-TextValue 'E47' '  -3.43%  '
Set-TextValue 'E48' '  -4.85%  '
Set-TextValue 'D49' '1.553'
Set-TextValue 'E49' '  -2.22%  '
Set-TextValue 'D50' '36.91'
Set-TextValue 'E50' '  -0.21%  '
Set-TextValue 'D51' '63.12'
Set-TextValue 'E51' '  -4.64%  '
